# feat: add 2022-Q4 data
#
# 1. Insert a new "2022-Q4" sheet right after "总计" (before "2022-Q3"),
#    populated with the fund-holding breakdown for that quarter.
# 2. Prepend a matching summary row to the "总计" sheet and shift the
#    existing rows (and their index column) down by one.

$wb = $excel.ActiveWorkbook

function Set-TextCell($cell, [string]$val) {
    # Writing a numeric-looking string via .Value lets the COM layer coerce
    # it to a real number (and drop leading zeros / trailing zeros), so we
    # use the classic "leading apostrophe" trick to force text storage for
    # anything that looks like a number. Genuine text is left alone so we
    # don't pick up a spurious quote-prefix style on cells that don't need it.
    if ($val -match '^-?[0-9]+(\.[0-9]+)?$') {
        $cell.Value = "'" + $val
    } else {
        $cell.Value = $val
    }
}

# ---------------------------------------------------------------------
# Step 1: build the new "2022-Q4" sheet by cloning "2022-Q3" (this keeps
# the exact same header/style/column layout) and overwriting its data.
# ---------------------------------------------------------------------

$srcSheet = $wb.Worksheets.Item("2022-Q3")
$srcSheet.Copy($srcSheet, $null)
$q4 = $wb.Worksheets.Item("2022-Q3 (2)")
$q4.Name = "2022-Q4"

$fundData = @(
    @("160921","大成多策略混合（LOF）A","8.25","87.54","7.59","0.6262","3"),
    @("016062","大成多策略混合（LOF）C","6.19","87.54","7.59","0.4698","3"),
    @("015208","信澳健康中国灵活配置混合C","12.04","93.44","3.54","0.4262","8"),
    @("003291","信澳健康中国灵活配置混合A","11.20","93.44","3.54","0.3965","8"),
    @("011598","信澳医药健康混合","7.67","93.37","3.44","0.2638","10"),
    @("200006","长城消费增值混合","5.59","92.96","4.53","0.2532","6"),
    @("000601","华宝创新优选混合","10.34","87.48","2.00","0.2068","9"),
    @("013037","长城大健康混合A","5.51","83.82","2.37","0.1306","9"),
    @("000523","国投瑞银医疗保健混合A","2.11","94.20","4.94","0.1042","4"),
    @("506008","长城科创两年定开混合A","3.19","77.75","2.70","0.0861","7"),
    @("014121","大成品质医疗股票A","0.56","89.89","6.94","0.0389","5"),
    @("010799","长城优选稳进六个月持有期混合A","1.96","32.14","1.57","0.0308","2"),
    @("000827","广发中证百度百发策略100指数E","2.68","92.98","0.99","0.0265","10"),
    @("005520","国投瑞银创新医疗混合","0.49","92.50","3.29","0.0161","10"),
    @("002681","金鹰元和灵活配置混合A","0.30","81.19","5.37","0.0161","4"),
    @("002682","金鹰元和灵活配置混合C","0.23","81.19","5.37","0.0124","4"),
    @("000649","长城久鑫灵活配置混合A","0.46","90.33","2.01","0.0092","10"),
    @("006890","上投摩根领先优选混合A","0.29","81.41","3.16","0.0092","6"),
    @("000826","广发中证百度百发策略100指数A","0.87","92.98","0.99","0.0086","10"),
    @("001318","东方新策略灵活配置混合A","0.39","36.37","2.11","0.0082","1"),
    @("013072","泰信医疗服务混合A","0.19","89.92","4.06","0.0077","7"),
    @("011082","国投瑞银医疗保健混合C","0.15","94.20","4.94","0.0074","4"),
    @("013038","长城大健康混合C","0.26","83.82","2.37","0.0062","9"),
    @("013073","泰信医疗服务混合C","0.15","89.92","4.06","0.0061","7"),
    @("014122","大成品质医疗股票C","0.08","89.89","6.94","0.0056","5"),
    @("400020","东方成长回报平衡混合","0.15","45.42","2.58","0.0039","1"),
    @("012793","长城科创两年定开混合C","0.11","77.75","2.70","0.0030","7"),
    @("002060","东方新策略灵活配置混合C","0.03","36.37","2.11","0.0006","1"),
    @("010800","长城优选稳进六个月持有期混合C","0.04","32.14","1.57","0.0006","2"),
    @("017461","长城久鑫灵活配置混合C","0.00","90.33","2.01","0","10"),
    @("017098","上投摩根领先优选混合C","0.00","81.41","3.16","0","6")
)

$rowCount = $fundData.Count      # 31 data rows -> rows 2..32
$lastRow = 1 + $rowCount

# The cloned sheet only has rows 1..25 (24 data rows). Extend formatting for
# the extra rows by copying the last existing data row's format downward.
if ($lastRow -gt 25) {
    $q4.Range("A25:H25").Copy()
    $q4.Range("A26:H$lastRow").PasteSpecial(-4122)
}

for ($i = 0; $i -lt $rowCount; $i++) {
    $r = 2 + $i
    $row = $fundData[$i]

    $q4.Cells.Item($r, 1).Value = $i            # A: 0-based index
    Set-TextCell $q4.Cells.Item($r, 2) $row[0]  # B: 基金代码
    Set-TextCell $q4.Cells.Item($r, 3) $row[1]  # C: 基金名称
    Set-TextCell $q4.Cells.Item($r, 4) $row[2]  # D: 基金规模
    Set-TextCell $q4.Cells.Item($r, 5) $row[3]  # E: 股票总仓位
    Set-TextCell $q4.Cells.Item($r, 6) $row[4]  # F: 仓位占比

    # G: 持有市值(亿元) - text everywhere except the two zero-holding funds
    # at the bottom of the sheet, which store a genuine numeric 0.
    if ($row[5] -eq "0") {
        $q4.Cells.Item($r, 7).Value = 0
    } else {
        Set-TextCell $q4.Cells.Item($r, 7) $row[5]
    }

    $q4.Cells.Item($r, 8).Value = [double]$row[6]  # H: 仓位排名 (number)
}

# ---------------------------------------------------------------------
# Step 2: add the 2022-Q4 summary row to the "总计" sheet, shifting the
# existing rows (and their 0-based index column) down by one.
# ---------------------------------------------------------------------

$total = $wb.Worksheets.Item("总计")

# Extend formatting for the new last row (row 6) by copying row 5's format.
$total.Range("A5:D5").Copy()
$total.Range("A6:D6").PasteSpecial(-4122)

$total.Cells.Item(6,1).Value = 4
$total.Cells.Item(6,2).Value = "2021-Q4"
$total.Cells.Item(6,3).Value = 2
$total.Cells.Item(6,4).Value = 0.02

$total.Cells.Item(5,1).Value = 3
$total.Cells.Item(5,2).Value = "2022-Q1"
$total.Cells.Item(5,3).Value = 3
$total.Cells.Item(5,4).Value = 0.08

$total.Cells.Item(4,1).Value = 2
$total.Cells.Item(4,2).Value = "2022-Q2"
$total.Cells.Item(4,3).Value = 14
$total.Cells.Item(4,4).Value = 0.87

$total.Cells.Item(3,1).Value = 1
$total.Cells.Item(3,2).Value = "2022-Q3"
$total.Cells.Item(3,3).Value = 24
$total.Cells.Item(3,4).Value = 1.67

$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q4"
$total.Cells.Item(2,3).Value = 31
$total.Cells.Item(2,4).Value = 3.18
